$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row heights
$ws.Rows.Item(2).RowHeight = 49.5
$ws.Rows.Item(3).RowHeight = 61.5
$ws.Rows.Item(4).RowHeight = 49.5
$ws.Rows.Item(5).RowHeight = 49.5
$ws.Rows.Item(6).RowHeight = 49.5
$ws.Rows.Item(7).RowHeight = 49.5
$ws.Rows.Item(8).RowHeight = 49.5
$ws.Rows.Item(9).RowHeight = 49.5
$ws.Rows.Item(10).RowHeight = 49.5

# Update selection
$ws.Range("L5").Select()
